# Update of financial figures on the WILC Income Statement / Balance Sheet /
# Cash Flow Statement worksheet (commit: "Doing Updates for Financials").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 86100
$ws.Range("E8").Value = 81200
$ws.Range("F8").Value = 86200
$ws.Range("G8").Value = 90700
$ws.Range("H8").Value = 92700
$ws.Range("I8").Value = 79000
$ws.Range("J8").Value = 72900
$ws.Range("D9").Value = 65600
$ws.Range("E9").Value = 60000
$ws.Range("F9").Value = 65500
$ws.Range("G9").Value = 68700
$ws.Range("H9").Value = 69600
$ws.Range("I9").Value = 60000
$ws.Range("J9").Value = 55900
$ws.Range("D10").Value = 20500
$ws.Range("E10").Value = 21100
$ws.Range("F10").Value = 20700
$ws.Range("G10").Value = 22000
$ws.Range("H10").Value = 23100
$ws.Range("I10").Value = 19000
$ws.Range("J10").Value = 17000
$ws.Range("D17").Value = 81400
$ws.Range("E17").Value = 74900
$ws.Range("F17").Value = 84300
$ws.Range("G17").Value = 84200
$ws.Range("H17").Value = 84700
$ws.Range("I17").Value = 72600
$ws.Range("J17").Value = 68200
$ws.Range("D18").Value = 4600
$ws.Range("E18").Value = 6300
$ws.Range("F18").Value = 1900
$ws.Range("G18").Value = 6500
$ws.Range("H18").Value = 8100
$ws.Range("I18").Value = 6500
$ws.Range("J18").Value = 4700
$ws.Range("D20").Value = 3900
$ws.Range("H20").Value = 3300
$ws.Range("D21").Value = 9500
$ws.Range("E21").Value = 5500
$ws.Range("F21").Value = 3600
$ws.Range("G21").Value = 8200
$ws.Range("H21").Value = 12600
$ws.Range("I21").Value = 9600
$ws.Range("J21").Value = 6000
$ws.Range("D23").Value = 8500
$ws.Range("G23").Value = 7200
$ws.Range("H23").Value = 11400
$ws.Range("I23").Value = 8800
$ws.Range("J23").Value = 5000
$ws.Range("D24").Value = 1600
$ws.Range("H24").Value = 2600
$ws.Range("I24").Value = 2100
$ws.Range("D26").Value = 6900
$ws.Range("G26").Value = 5200
$ws.Range("H26").Value = 8800
$ws.Range("I26").Value = 6600
$ws.Range("D27").Value = 6900
$ws.Range("G27").Value = 5200
$ws.Range("H27").Value = 8800
$ws.Range("I27").Value = 6600
$ws.Range("D32").Value = -3900
$ws.Range("H32").Value = -3300
$ws.Range("D33").Value = 6900
$ws.Range("G33").Value = 5200
$ws.Range("H33").Value = 8800
$ws.Range("I33").Value = 6600
$ws.Range("D35").Value = 6900
$ws.Range("G35").Value = 5200
$ws.Range("H35").Value = 8800
$ws.Range("I35").Value = 6600
$ws.Range("E41").Value = 24900
$ws.Range("F41").Value = 5500
$ws.Range("G41").Value = 7700
$ws.Range("I41").Value = 5900
$ws.Range("D42").Value = 67200
$ws.Range("E42").Value = 39700
$ws.Range("F42").Value = 62000
$ws.Range("G42").Value = 54400
$ws.Range("H42").Value = 38400
$ws.Range("I42").Value = 53800
$ws.Range("J42").Value = 49700
$ws.Range("D43").Value = 25200
$ws.Range("E43").Value = 23300
$ws.Range("F43").Value = 24700
$ws.Range("G43").Value = 24800
$ws.Range("H43").Value = 41400
$ws.Range("I43").Value = 21200
$ws.Range("J43").Value = 20100
$ws.Range("D44").Value = 11000
$ws.Range("E44").Value = 11600
$ws.Range("F44").Value = 9500
$ws.Range("G44").Value = 13400
$ws.Range("H44").Value = 14900
$ws.Range("I44").Value = 13600
$ws.Range("J44").Value = 9000
$ws.Range("D46").Value = 109000
$ws.Range("E46").Value = 101200
$ws.Range("F46").Value = 102300
$ws.Range("G46").Value = 100800
$ws.Range("H46").Value = 97700
$ws.Range("I46").Value = 94600
$ws.Range("J46").Value = 83900
$ws.Range("D48").Value = 11400
$ws.Range("E48").Value = 11700
$ws.Range("F48").Value = 12200
$ws.Range("G48").Value = 12500
$ws.Range("H48").Value = 11300
$ws.Range("I48").Value = 11500
$ws.Range("J48").Value = 11700
$ws.Range("E52").Value = 600
$ws.Range("D54").Value = 120500
$ws.Range("E54").Value = 113500
$ws.Range("F54").Value = 115600
$ws.Range("G54").Value = 113500
$ws.Range("H54").Value = 109000
$ws.Range("I54").Value = 106100
$ws.Range("J54").Value = 95900
$ws.Range("D57").Value = 3400
$ws.Range("H57").Value = 5700
$ws.Range("I57").Value = 8300
$ws.Range("J57").Value = 7300
$ws.Range("I58").Value = 2700
$ws.Range("D59").Value = 2100
$ws.Range("G59").Value = 2200
$ws.Range("E60").Value = 5400
$ws.Range("F60").Value = 5100
$ws.Range("G60").Value = 6800
$ws.Range("H60").Value = 7700
$ws.Range("I60").Value = 13900
$ws.Range("J60").Value = 10200
$ws.Range("D66").Value = 5900
$ws.Range("E66").Value = 5600
$ws.Range("G66").Value = 7000
$ws.Range("H66").Value = 8100
$ws.Range("I66").Value = 14100
$ws.Range("J66").Value = 10300
$ws.Range("D72").Value = 79000
$ws.Range("E72").Value = 72100
$ws.Range("F72").Value = 74500
$ws.Range("G72").Value = 72600
$ws.Range("H72").Value = 67400
$ws.Range("I72").Value = 58600
$ws.Range("J72").Value = 52000
$ws.Range("D76").Value = 114700
$ws.Range("E76").Value = 107900
$ws.Range("F76").Value = 110300
$ws.Range("G76").Value = 106500
$ws.Range("H76").Value = 100900
$ws.Range("I76").Value = 92100
$ws.Range("J76").Value = 85600
$ws.Range("D81").Value = 6900
$ws.Range("G81").Value = 5200
$ws.Range("H81").Value = 8800
$ws.Range("I81").Value = 6600
$ws.Range("E83").Value = 1000
$ws.Range("G89").Value = 5400
$ws.Range("J89").Value = 10100
$ws.Range("D94").Value = -8500
$ws.Range("E94").Value = 14400
$ws.Range("F94").Value = -6800
$ws.Range("G94").Value = 7500
$ws.Range("I94").Value = 5700
$ws.Range("J94").Value = -28800
$ws.Range("E96").Value = -5300
$ws.Range("E100").Value = -5300
$ws.Range("H100").Value = -2700
$ws.Range("I100").Value = 2500
$ws.Range("E102").Value = 13800
$ws.Range("G102").Value = 12900
$ws.Range("H102").Value = -5900
$ws.Range("I102").Value = 6300
$ws.Range("J102").Value = -21800
